# Commit 10: criado a view com o detalhe dos atrasos no relatório
#
# Adds a new student ("Aluno 21") row to the "teste_alunos" sheet and
# updates the current selection to cover the newly added data, matching
# the row layout already used by the previous entries (A:R.A.,
# B:Nome do estudante, C:Série/turma, D:Endereço, E:Responsável 1,
# F:Responsável 2, G:Contato(s)).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 202521
$ws.Range("B12").Value = "Aluno 21"
$ws.Range("C12").Value = "7B"
$ws.Range("D12").Value = "Rua do Aluno 21"
$ws.Range("E12").Value = "Pai do Aluno 21"
$ws.Range("F12").Value = "Mãe do Aluno 21"
$ws.Range("G12").Value = 11987654341

# Mirror the selection left behind in the saved workbook (the author
# selected the newly added rows/range before saving).
$ws.Range("A9:G12").Select()

$wb.Save()
